# feat: add 2022-Q1 data
# - Insert a new worksheet "2022-Q1" (holding-detail table) right before the
#   "总计" (totals) summary sheet.
# - Populate the new sheet with the 2022-Q1 fund holding detail row.
# - Prepend a new "2022-Q1" summary row into the "总计" sheet, shifting the
#   existing history rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right before "总计".
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2021-Q4")

$newSheet = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$newSheet.Name = "2022-Q1"

# NOTE: re-fetch "总计" by name after Add() — the sheet object/variable that
# was passed in as the "Before" argument above gets repointed to the newly
# inserted sheet, so it can no longer be relied on to address "总计".
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header formatting (bold / border / centered) + the index-column style
# from the previous quarter sheet, which has the identical B:H layout.
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row text
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Data row (2022-Q1 top holding). Columns B/D/E/F/G are numeric-looking text
# in the source data (matches sibling quarter sheets), so we force the
# "Text" number format while assigning, then reset the cell style back to
# Normal so no residual style index is left attached to the cell.
$newSheet.Cells.Item(2, 1).Value = 0

$newSheet.Cells.Item(2, 2).NumberFormat = "@"
$newSheet.Cells.Item(2, 2).Value = "206009"
$newSheet.Cells.Item(2, 2).Style = "Normal"

$newSheet.Cells.Item(2, 3).Value = "鹏华新兴产业混合"

$newSheet.Cells.Item(2, 4).NumberFormat = "@"
$newSheet.Cells.Item(2, 4).Value = "44.95"
$newSheet.Cells.Item(2, 4).Style = "Normal"

$newSheet.Cells.Item(2, 5).NumberFormat = "@"
$newSheet.Cells.Item(2, 5).Value = "90.17"
$newSheet.Cells.Item(2, 5).Style = "Normal"

$newSheet.Cells.Item(2, 6).NumberFormat = "@"
$newSheet.Cells.Item(2, 6).Value = "5.45"
$newSheet.Cells.Item(2, 6).Style = "Normal"

$newSheet.Cells.Item(2, 7).NumberFormat = "@"
$newSheet.Cells.Item(2, 7).Value = "2.4498"
$newSheet.Cells.Item(2, 7).Style = "Normal"

$newSheet.Cells.Item(2, 8).Value = 4

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: add a new 2022-Q1 row on top, push history down.
# ---------------------------------------------------------------------------
$ws = $totalSheet

# Extend the numbered index column (A) down by one row, copying its style.
$ws.Cells.Item(6, 1).Copy()
$ws.Cells.Item(7, 1).PasteSpecial(-4122)
$ws.Cells.Item(7, 1).Value = 5

# Shift the quarter / count / value columns down by one row (bottom-up).
$ws.Cells.Item(7, 2).Value = "2020-Q4"
$ws.Cells.Item(7, 3).Value = 21
$ws.Cells.Item(7, 4).Value = 6.57

$ws.Cells.Item(6, 2).Value = "2021-Q1"
$ws.Cells.Item(6, 3).Value = 17
$ws.Cells.Item(6, 4).Value = 8.970000000000001

$ws.Cells.Item(5, 2).Value = "2021-Q2"
$ws.Cells.Item(5, 3).Value = 12
$ws.Cells.Item(5, 4).Value = 6.96

$ws.Cells.Item(4, 2).Value = "2021-Q3"
$ws.Cells.Item(4, 3).Value = 11
$ws.Cells.Item(4, 4).Value = 6.35

$ws.Cells.Item(3, 2).Value = "2021-Q4"
$ws.Cells.Item(3, 3).Value = 10
$ws.Cells.Item(3, 4).Value = 5.46

$ws.Cells.Item(2, 2).Value = "2022-Q1"
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 2.45
